$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Tgfb1"
$ws.Cells.Item(2,3).Value = "Itgb8"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 35.906979
$ws.Cells.Item(2,8).Value = 107.720937
$ws.Cells.Item(2,9).Value = 0.6107087147789413
$ws.Cells.Item(2,10).Value = 0.6107087147789412
$ws.Cells.Item(2,11).Value = 1
$ws.Cells.Item(2,12).Value = 0.3333333333333333
$ws.Cells.Item(2,13).Value = 0.027767
$ws.Cells.Item(2,14).Value = 0.083301
$ws.Cells.Item(2,15).Value = 0.002923627791763407
$ws.Cells.Item(2,16).Value = 0.002923627791763407
$ws.Cells.Item(2,17).Value = 0.997029085893
$ws.Cells.Item(2,18).Value = 8.973261773036999
$ws.Cells.Item(2,19).Value = 0.001785484971199825
$ws.Cells.Item(2,20).Value = 0.001785484971199824

# Row 3
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Tgfb1"
$ws.Cells.Item(3,3).Value = "Itgb8"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 35.906979
$ws.Cells.Item(3,8).Value = 107.720937
$ws.Cells.Item(3,9).Value = 0.6107087147789413
$ws.Cells.Item(3,10).Value = 0.6107087147789412
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 4.237840333333334
$ws.Cells.Item(3,14).Value = 12.713521
$ws.Cells.Item(3,15).Value = 0.4462083687682946
$ws.Cells.Item(3,16).Value = 0.4462083687682946
$ws.Cells.Item(3,17).Value = 152.168043854353
$ws.Cells.Item(3,18).Value = 1369.512394689177
$ws.Cells.Item(3,19).Value = 0.2725033394140931
$ws.Cells.Item(3,20).Value = 0.272503339414093

# Row 4
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Tgfb1"
$ws.Cells.Item(4,3).Value = "Itgb8"
$ws.Cells.Item(4,4).Value = "MuSCs"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 35.906979
$ws.Cells.Item(4,8).Value = 107.720937
$ws.Cells.Item(4,9).Value = 0.6107087147789413
$ws.Cells.Item(4,10).Value = 0.6107087147789412
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 5.231839666666667
$ws.Cells.Item(4,14).Value = 15.695519
$ws.Cells.Item(4,15).Value = 0.5508680034399419
$ws.Cells.Item(4,16).Value = 0.5508680034399419
$ws.Cells.Item(4,17).Value = 187.859557042367
$ws.Cells.Item(4,18).Value = 1690.736013381303
$ws.Cells.Item(4,19).Value = 0.3364198903936483
$ws.Cells.Item(4,20).Value = 0.3364198903936483

# Row 5
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Tgfb1"
$ws.Cells.Item(5,3).Value = "Itgb8"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 17.04862266666667
$ws.Cells.Item(5,8).Value = 51.14586800000001
$ws.Cells.Item(5,9).Value = 0.2899643113254147
$ws.Cells.Item(5,10).Value = 0.2899643113254147
$ws.Cells.Item(5,11).Value = 1
$ws.Cells.Item(5,12).Value = 0.3333333333333333
$ws.Cells.Item(5,13).Value = 0.027767
$ws.Cells.Item(5,14).Value = 0.083301
$ws.Cells.Item(5,15).Value = 0.002923627791763407
$ws.Cells.Item(5,16).Value = 0.002923627791763407
$ws.Cells.Item(5,17).Value = 0.4733891055853334
$ws.Cells.Item(5,18).Value = 4.260501950268001
$ws.Cells.Item(5,19).Value = 0.0008477477192105194
$ws.Cells.Item(5,20).Value = 0.0008477477192105192

# Row 6
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Tgfb1"
$ws.Cells.Item(6,3).Value = "Itgb8"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 17.04862266666667
$ws.Cells.Item(6,8).Value = 51.14586800000001
$ws.Cells.Item(6,9).Value = 0.2899643113254147
$ws.Cells.Item(6,10).Value = 0.2899643113254147
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 4.237840333333334
$ws.Cells.Item(6,14).Value = 12.713521
$ws.Cells.Item(6,15).Value = 0.4462083687682946
$ws.Cells.Item(6,16).Value = 0.4462083687682946
$ws.Cells.Item(6,17).Value = 72.24934076458091
$ws.Cells.Item(6,18).Value = 650.2440668812282
$ws.Cells.Item(6,19).Value = 0.1293845023575352
$ws.Cells.Item(6,20).Value = 0.1293845023575352

# Row 7
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Tgfb1"
$ws.Cells.Item(7,3).Value = "Itgb8"
$ws.Cells.Item(7,4).Value = "MuSCs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 17.04862266666667
$ws.Cells.Item(7,8).Value = 51.14586800000001
$ws.Cells.Item(7,9).Value = 0.2899643113254147
$ws.Cells.Item(7,10).Value = 0.2899643113254147
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 5.231839666666667
$ws.Cells.Item(7,14).Value = 15.695519
$ws.Cells.Item(7,15).Value = 0.5508680034399419
$ws.Cells.Item(7,16).Value = 0.5508680034399419
$ws.Cells.Item(7,17).Value = 89.19566032949913
$ws.Cells.Item(7,18).Value = 802.7609429654922
$ws.Cells.Item(7,19).Value = 0.159732061248669
$ws.Cells.Item(7,20).Value = 0.1597320612486689

# Row 8
$ws.Cells.Item(8,1).Value = "MuSCs"
$ws.Cells.Item(8,2).Value = "Tgfb1"
$ws.Cells.Item(8,3).Value = "Itgb8"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 5.839988000000001
$ws.Cells.Item(8,8).Value = 17.519964
$ws.Cells.Item(8,9).Value = 0.09932697389564409
$ws.Cells.Item(8,10).Value = 0.09932697389564407
$ws.Cells.Item(8,11).Value = 1
$ws.Cells.Item(8,12).Value = 0.3333333333333333
$ws.Cells.Item(8,13).Value = 0.027767
$ws.Cells.Item(8,14).Value = 0.083301
$ws.Cells.Item(8,15).Value = 0.002923627791763407
$ws.Cells.Item(8,16).Value = 0.002923627791763407
$ws.Cells.Item(8,17).Value = 0.162158946796
$ws.Cells.Item(8,18).Value = 1.459430521164
$ws.Cells.Item(8,19).Value = 0.0002903951013530635
$ws.Cells.Item(8,20).Value = 0.0002903951013530635

# Row 9
$ws.Cells.Item(9,1).Value = "MuSCs"
$ws.Cells.Item(9,2).Value = "Tgfb1"
$ws.Cells.Item(9,3).Value = "Itgb8"
$ws.Cells.Item(9,4).Value = "FAPs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 5.839988000000001
$ws.Cells.Item(9,8).Value = 17.519964
$ws.Cells.Item(9,9).Value = 0.09932697389564409
$ws.Cells.Item(9,10).Value = 0.09932697389564407
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 4.237840333333334
$ws.Cells.Item(9,14).Value = 12.713521
$ws.Cells.Item(9,15).Value = 0.4462083687682946
$ws.Cells.Item(9,16).Value = 0.4462083687682946
$ws.Cells.Item(9,17).Value = 24.74893669258267
$ws.Cells.Item(9,18).Value = 222.740430233244
$ws.Cells.Item(9,19).Value = 0.04432052699666633
$ws.Cells.Item(9,20).Value = 0.04432052699666632

# Row 10
$ws.Cells.Item(10,1).Value = "MuSCs"
$ws.Cells.Item(10,2).Value = "Tgfb1"
$ws.Cells.Item(10,3).Value = "Itgb8"
$ws.Cells.Item(10,4).Value = "MuSCs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 5.839988000000001
$ws.Cells.Item(10,8).Value = 17.519964
$ws.Cells.Item(10,9).Value = 0.09932697389564409
$ws.Cells.Item(10,10).Value = 0.09932697389564407
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 5.231839666666667
$ws.Cells.Item(10,14).Value = 15.695519
$ws.Cells.Item(10,15).Value = 0.5508680034399419
$ws.Cells.Item(10,16).Value = 0.5508680034399419
$ws.Cells.Item(10,17).Value = 30.55388087125734
$ws.Cells.Item(10,18).Value = 274.9849278413161
$ws.Cells.Item(10,19).Value = 0.05471605179762468
$ws.Cells.Item(10,20).Value = 0.05471605179762468
